# Fix Training Data Issue (#48)
# The "Date" column (BF) was off by one day because of how NBA stats were
# reported. Correct every data row's date string from the old
# "6-25-2013-14" placeholder to the proper ISO date "2014-06-25".
# Row 1 is the "Date" header and stays untouched; rows 2-31 hold the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 31
$col      = 58   # column BF

$correctedDate = "2014-06-25"

$rng = $ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col))

# Pre-format as Text so Excel doesn't silently reinterpret the ISO-looking
# string "2014-06-25" as a date serial number while we type it in.
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = $correctedDate
}

# The cells never had any special formatting of their own (no border/fill/
# alignment), so drop the temporary Text format again to leave the cells
# styled exactly as they were before - just with the corrected value.
$rng.ClearFormats()
